$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1035065.75
$ws.Range("I17").Value = 963
$ws.Range("J17").Value = 1051223.6
$ws.Range("K17").Value = 2889
$ws.Range("L17").Value = 3153670.8
$ws.Range("M17").Value = -2721
$ws.Range("N17").Value = -3154006.8
$ws.Range("H113").Value = 4975.5625
$ws.Range("I113").Value = 4698
$ws.Range("K113").Value = 4698
$ws.Range("M113").Value = -1444
$ws.Range("H125").Value = 4505.6665
$ws.Range("I125").Value = 7182
$ws.Range("K125").Value = 64638
$ws.Range("M125").Value = -62178
$ws.Range("H137").Value = 19065.75
$ws.Range("I137").Value = 19544.21
$ws.Range("J137").Value = 17247.6
$ws.Range("K137").Value = 58632.63
$ws.Range("L137").Value = 51742.8
$ws.Range("M137").Value = -56082.63
$ws.Range("N137").Value = -56842.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4284.5
$ws.Range("I45").Value = 3464.3333
$ws.Range("K45").Value = 3464.3333
$ws.Range("M45").Value = -3087.3333
$ws.Range("H61").Value = 3342.0386
$ws.Range("I61").Value = 3091.875
$ws.Range("J61").Value = 3742.3
$ws.Range("K61").Value = 3091.875
$ws.Range("L61").Value = 3742.3
$ws.Range("M61").Value = -2879.875
$ws.Range("N61").Value = -4166.3
$ws.Range("H74").Value = 34630.582
$ws.Range("I74").Value = 43779.145
$ws.Range("K74").Value = 43779.145
$ws.Range("M74").Value = -42905.145
$ws.Range("H77").Value = 34630.582
$ws.Range("I77").Value = 43779.145
$ws.Range("K77").Value = 218895.725
$ws.Range("M77").Value = -214527.725
$ws.Range("H122").Value = 2715.423
$ws.Range("I122").Value = 2459.1365
$ws.Range("K122").Value = 7377.4095
$ws.Range("M122").Value = -4927.4095
$ws.Range("H135").Value = 115897.125
$ws.Range("J135").Value = 115897.125
$ws.Range("L135").Value = 115897.125
$ws.Range("N135").Value = -126037.125
$ws.Range("H136").Value = 3342.0386
$ws.Range("I136").Value = 3091.875
$ws.Range("J136").Value = 3742.3
$ws.Range("K136").Value = 9275.625
$ws.Range("L136").Value = 11226.9
$ws.Range("M136").Value = -6725.625
$ws.Range("N136").Value = -16326.9
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1165.2142
$ws.Range("I94").Value = 574
$ws.Range("K94").Value = 574
$ws.Range("M94").Value = -123
$ws.Range("H99").Value = 3962
$ws.Range("I99").Value = 2336.6667
$ws.Range("K99").Value = 2336.6667
$ws.Range("M99").Value = -838.6667000000002
$ws.Range("H105").Value = 8522.546
$ws.Range("I105").Value = 9202.632
$ws.Range("J105").Value = 7599.5713
$ws.Range("K105").Value = 9202.632
$ws.Range("L105").Value = 7599.5713
$ws.Range("M105").Value = -7455.632
$ws.Range("N105").Value = -11093.5713
$ws.Range("H107").Value = 3345.875
$ws.Range("I107").Value = 3118.4614
$ws.Range("K107").Value = 3118.4614
$ws.Range("M107").Value = -1198.4614
$ws.Range("H113").Value = 3974.5
$ws.Range("I113").Value = 3974.5
$ws.Range("K113").Value = 3974.5
$ws.Range("M113").Value = -1804.5
$ws.Range("H132").Value = 139780
$ws.Range("J132").Value = 139780
$ws.Range("L132").Value = 139780
$ws.Range("N132").Value = -149900
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 26833.334
$ws.Range("J50").Value = 34333.332
$ws.Range("L50").Value = 34333.332
$ws.Range("N50").Value = -35583.332
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 42.8125
$ws.Range("I2").Value = 51.8
$ws.Range("J2").Value = 27.833334
$ws.Range("K2").Value = 310.8
$ws.Range("L2").Value = 167.000004
$ws.Range("M2").Value = -197.8
$ws.Range("N2").Value = -393.000004
$ws.Range("H4").Value = 67426150
$ws.Range("I4").Value = 101772750
$ws.Range("J4").Value = 10181830
$ws.Range("K4").Value = 305318250
$ws.Range("L4").Value = 30545490
$ws.Range("M4").Value = -305318138
$ws.Range("N4").Value = -30545714
$ws.Range("H12").Value = 1074.9
$ws.Range("I12").Value = 400
$ws.Range("J12").Value = 1364.1428
$ws.Range("K12").Value = 1200
$ws.Range("L12").Value = 4092.4284
$ws.Range("M12").Value = -1027
$ws.Range("N12").Value = -4438.428400000001
$ws.Range("H18").Value = 4247.5
$ws.Range("I18").Value = 2330
$ws.Range("K18").Value = 6990
$ws.Range("M18").Value = -6821
$ws.Range("H34").Value = 616.7857
$ws.Range("I34").Value = 130
$ws.Range("J34").Value = 1265.8334
$ws.Range("K34").Value = 390
$ws.Range("L34").Value = 3797.5002
$ws.Range("M34").Value = -306
$ws.Range("N34").Value = -3965.5002
$ws.Range("H55").Value = 1127924.8
$ws.Range("J55").Value = 4599.6
$ws.Range("L55").Value = 13798.8
$ws.Range("N55").Value = -14152.8
$ws.Range("H102").Value = 5249.5
$ws.Range("J102").Value = 7500
$ws.Range("L102").Value = 22500
$ws.Range("N102").Value = -27368
$ws.Range("H107").Value = 936.2727
$ws.Range("J107").Value = 1493
$ws.Range("L107").Value = 4479
$ws.Range("N107").Value = -8319
$ws.Range("H117").Value = 2403.5
$ws.Range("I117").Value = 871.3333
$ws.Range("J117").Value = 7000
$ws.Range("K117").Value = 2613.9999
$ws.Range("L117").Value = 21000
$ws.Range("M117").Value = 828.0001000000002
$ws.Range("N117").Value = -27884
$ws.Range("H121").Value = 383.33334
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H125").Value = 6875
$ws.Range("I125").Value = 6875
$ws.Range("K125").Value = 20625
$ws.Range("M125").Value = -15705
$ws.Range("H129").Value = 2002.7097
$ws.Range("J129").Value = 2654.9473
$ws.Range("L129").Value = 7964.841899999999
$ws.Range("N129").Value = -17964.8419
$ws.Range("H130").Value = 2699.2
$ws.Range("J130").Value = 3116.5
$ws.Range("L130").Value = 9349.5
$ws.Range("N130").Value = -19389.5
$ws.Range("H131").Value = 1658.3572
$ws.Range("J131").Value = 1700.2727
$ws.Range("L131").Value = 5100.8181
$ws.Range("N131").Value = -15180.8181
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1020
$ws.Range("J97").Value = 1262.5555
$ws.Range("L97").Value = 1262.5555
$ws.Range("N97").Value = -2254.5555
$ws.Range("H107").Value = 589.1818
$ws.Range("I107").Value = 579.5714
$ws.Range("J107").Value = 606
$ws.Range("K107").Value = 579.5714
$ws.Range("L107").Value = 606
$ws.Range("M107").Value = 1340.4286
$ws.Range("N107").Value = -4446
$ws.Range("H122").Value = 1968.2
$ws.Range("I122").Value = 1272.1875
$ws.Range("J122").Value = 3205.5557
$ws.Range("K122").Value = 3816.5625
$ws.Range("L122").Value = 9616.667099999999
$ws.Range("M122").Value = -1366.5625
$ws.Range("N122").Value = -14516.6671
$ws.Range("H139").Value = 110000
$ws.Range("J139").Value = 110000
$ws.Range("L139").Value = 110000
$ws.Range("N139").Value = -120280
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 822.8889
$ws.Range("I16").Value = 773.4545000000001
$ws.Range("K16").Value = 773.4545000000001
$ws.Range("M16").Value = -603.4545000000001
$ws.Range("H20").Value = 769.5
$ws.Range("I20").Value = 305
$ws.Range("K20").Value = 305
$ws.Range("M20").Value = -79
$ws.Range("H22").Value = 1513.1818
$ws.Range("I22").Value = 1112.2106
$ws.Range("K22").Value = 1112.2106
$ws.Range("M22").Value = -817.2106000000001
$ws.Range("H27").Value = 1513.1818
$ws.Range("I27").Value = 1112.2106
$ws.Range("K27").Value = 1112.2106
$ws.Range("M27").Value = -1005.2106
$ws.Range("H43").Value = 22638.889
$ws.Range("J43").Value = 22821.428
$ws.Range("L43").Value = 22821.428
$ws.Range("N43").Value = -23207.428
$ws.Range("H122").Value = 3810.647
$ws.Range("I122").Value = 2820.4285
$ws.Range("J122").Value = 8431.666999999999
$ws.Range("K122").Value = 8461.2855
$ws.Range("L122").Value = 25295.001
$ws.Range("M122").Value = -6011.2855
$ws.Range("N122").Value = -30195.001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3343
$ws.Range("I96").Value = 2993
$ws.Range("J96").Value = 3722.1667
$ws.Range("K96").Value = 2993
$ws.Range("L96").Value = 3722.1667
$ws.Range("M96").Value = -1620
$ws.Range("N96").Value = -6468.1667
$ws.Range("H136").Value = 3107.3115
$ws.Range("I136").Value = 2246.689
$ws.Range("J136").Value = 5527.8125
$ws.Range("K136").Value = 6740.066999999999
$ws.Range("L136").Value = 16583.4375
$ws.Range("M136").Value = -4190.066999999999
$ws.Range("N136").Value = -21683.4375

Write-Host "Applied all Famfrit_Profits market data updates"